$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Team")

# Update column G (totaltimetaken) values from 20 to 30 for rows 2-11,
# reflecting the corrected sample data per the ER diagram.
for ($row = 2; $row -le 11; $row++) {
    $ws.Cells.Item($row, 7).Value = 30
}

# Move the selection from F11 to G11.
$ws.Range("G11").Select()
